# Apply translation updates to the Tehilim (Psalm 46) data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of old translation text -> new translation text for the "translation" column (B).
$replacements = @{
    "Начальнику хора" = "Руководителю хора"
    "Кореевых"        = "Кораха"
    "на девятых"      = "на Аламот"
    "найденная"       = "доступный"
    "весьма"          = "очень"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = $ws.UsedRange.Rows.Count }

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($null -ne $val -and $replacements.ContainsKey([string]$val)) {
        $cell.Value = $replacements[[string]$val]
    }
}
